$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settlements")

# Append the new settlement record as row 6 (below the existing 5 data rows).
$row = 6

$ws.Cells.Item($row, 1).Value  = "settlement_1753294282531_mmatj9jez"
$ws.Cells.Item($row, 2).Value  = "user_1753125931723_8ftkkx2pf"
$ws.Cells.Item($row, 3).Value  = "Adarsh"
$ws.Cells.Item($row, 4).Value  = 9.97
$ws.Cells.Item($row, 5).Value  = "lucky@okaxis"
$ws.Cells.Item($row, 6).Value  = "ORD-001,ORD-002,ORD-004"
$ws.Cells.Item($row, 7).Value  = 3
$ws.Cells.Item($row, 8).Value  = "INR"
$ws.Cells.Item($row, 9).Value  = "approved"
$ws.Cells.Item($row, 10).Value = "settled_partially"
$ws.Cells.Item($row, 11).Value = "2025-07-23T18:11:22.531Z"
$ws.Cells.Item($row, 12).Value = "2025-07-23T18:11:57.047Z"
$ws.Cells.Item($row, 13).Value = 5
# transactionId "111" looks numeric - force it to stay text (matches source t="str").
$ws.Cells.Item($row, 14).Value = "'111"
$ws.Cells.Item($row, 15).Value = "payment-proof-1753294317020-14169049.webp"
$ws.Cells.Item($row, 16).Value = "user_1753040616422_hgtapju6r"
$ws.Cells.Item($row, 17).Value = "2025-07-23T18:11:57.043Z"
